$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matrix Multiplication")
$ws.Activate()
